# Clean up the "Authors" column (column E) data: the author-list strings
# stored in this sheet use a comma-separated format where each separator
# should be followed by an extra space (i.e. ",            Name" becomes
# ",             Name" - one more space than before). This mirrors the
# author's commit "Cleaning more data. Fixed an issue where some files
# were displaying incorrect info." which re-saved every Authors cell with
# corrected spacing after each comma.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dims = $ws.UsedRange
$lastRow = $dims.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)   # column E = 5th column ("Authors")
    $current = $cell.Value2
    if ($current -ne $null -and $current.GetType().Name -eq "String" -and $current.Length -gt 0) {
        $updated = [System.Text.RegularExpressions.Regex]::Replace($current, ",(\s+)", ", `$1")
        if ($updated -ne $current) {
            $cell.Value = $updated
        }
    }
}
